$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '95.262.99'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '3.556.08'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '234.90'
$ws.Range('E5').Value = '  -2.38%  '
$ws.Range('D6').Value = '648.89'
$ws.Range('E6').Value = '  +1.71%  '
$ws.Range('E7').Value = '  -2.04%  '
$ws.Range('E8').Value = '  -1.90%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('D10').Value = '0.979'
$ws.Range('E10').Value = '  -4.31%  '
$ws.Range('D11').Value = '3.556.67'
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').Value = '41.96'
$ws.Range('E13').Value = '  -4.53%  '
$ws.Range('D14').Value = '6.46'
$ws.Range('E14').Value = '  +2.28%  '
$ws.Range('D15').Value = '4.218.06'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '95.185.18'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('D18').Value = '3.563.65'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').Value = '7.86'
$ws.Range('E19').Value = '  -6.72%  '
$ws.Range('D20').Value = '12.64'
$ws.Range('E20').Value = '  -2.67%  '
$ws.Range('D21').Value = '17.61'
$ws.Range('E21').Value = '  -3.38%  '
$ws.Range('D22').Value = '3.45'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').Value = '502.61'
$ws.Range('E23').Value = '  -2.81%  '
$ws.Range('D24').Value = '0.468'
$ws.Range('E24').Value = '  -8.75%  '
$ws.Range('E25').Value = '  -2.33%  '
$ws.Range('D26').Value = '6.53'
$ws.Range('E26').Value = '  -3.90%  '
$ws.Range('D27').Value = '94.22'
$ws.Range('E27').Value = '  -2.61%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '3.748.31'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = '12.37'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('D30').Value = '3.01'
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('D32').Value = '11.15'
$ws.Range('E32').Value = '  -4.22%  '
$ws.Range('E33').Value = '  -4.60%  '
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('E35').Value = '  -3.48%  '
$ws.Range('D36').Value = '31.60'
$ws.Range('E36').Value = '  +4.17%  '
$ws.Range('D37').Value = '0.552'
$ws.Range('E37').Value = '  -3.45%  '
$ws.Range('D38').Value = '8.13'
$ws.Range('E38').Value = '  +5.69%  '
$ws.Range('D39').Value = '558.52'
$ws.Range('E39').Value = '  -5.74%  '
$ws.Range('D40').Value = '1.50'
$ws.Range('E40').Value = '  +2.96%  '
$ws.Range('E42').Value = '  -1.93%  '
$ws.Range('D43').Value = '0.894'
$ws.Range('E43').Value = '  -4.37%  '
$ws.Range('D44').Value = '1.75'
$ws.Range('E44').Value = '  +0.40%  '
$ws.Range('D45').Value = '34.13'
$ws.Range('E45').Value = '  +32.42%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Value = '23.53'
$ws.Range('E46').Value = '  -1.31%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value = '2.25'
$ws.Range('E47').Value = '  +2.80%  '
$ws.Range('D48').Value = '5.58'
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('B49').Value = 'MantraDAO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D49').Value = '3.58'
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0406'
$ws.Range('E50').Value = '  -5.27%  '
$ws.Range('D51').Value = '53.29'
$ws.Range('E51').Value = '  -1.65%  '
